$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two new pieces of text that get introduced by this edit (they become
# new shared-string entries #29 and #30).
$textAvg  = "Calculate VAM by the average school test score"
$textPoly = "Add also only a third-order`npolynomial in the prior year’s reading and math`nscores"

# Row 4 ("1" -> descriptive text) for the three repeated blocks (B/G/L).
$ws.Range("B4").Value = $textAvg
$ws.Range("G4").Value = $textAvg
$ws.Range("L4").Value = $textAvg

# Row 5 ("1" -> descriptive text).
$ws.Range("B5").Value = $textAvg
$ws.Range("G5").Value = $textAvg
$ws.Range("L5").Value = $textAvg

# Row 6 ("1" -> descriptive text).
$ws.Range("B6").Value = $textAvg
$ws.Range("G6").Value = $textAvg
$ws.Range("L6").Value = $textAvg

# Row 7 ("2" -> long wrapped text).
$ws.Range("B7").Value = $textPoly
$ws.Range("G7").Value = $textPoly
$ws.Range("L7").Value = $textPoly

# Row 8 ("2" -> long wrapped text).
$ws.Range("B8").Value = $textPoly
$ws.Range("G8").Value = $textPoly
$ws.Range("L8").Value = $textPoly

# Row 9 ("2" -> long wrapped text).
$ws.Range("B9").Value = $textPoly
$ws.Range("G9").Value = $textPoly
$ws.Range("L9").Value = $textPoly

# The long text cells (rows 7-9) get a new wrap-text style (cellXfs index 3
# in the target) and a taller row to fit the wrapped paragraph.
$ws.Range("B7:B9").WrapText = $true
$ws.Range("G7:G9").WrapText = $true
$ws.Range("L7:L9").WrapText = $true

$ws.Rows(7).RowHeight = 156.75
$ws.Rows(8).RowHeight = 156.75
$ws.Rows(9).RowHeight = 156.75

# Column L (12) is widened to fit the new text (bestFit column in the diff).
$ws.Columns(12).ColumnWidth = 22.83333333333

# Selection moves from N2 to B8.
$ws.Range("B8").Select()
